# Fix parse numeric sku bug
#
# "sku" values that look numeric (e.g. "1234a" was mis-typed/mis-parsed as
# numeric) were being written into the sheet as text instead of a real
# number. Sheet1!A2 held the sku for Product 1 as the text string "1234a"
# even though the other sku in the sheet (A3, "1234b") is a correctly
# formed text sku. Looking at the data, A2's value was simply wrong - it
# should have been the plain number 1234.
#
# Correct it by writing the real numeric value into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 1234

# Reflect where the selection ended up after making the fix.
[void]$ws.Range("B7").Select()
